$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.851962457337884
$ws.Range("C2").Value = 0.8212616822429907
$ws.Range("D2").Value = 0.8997440273037542
$ws.Range("E2").Value = 0.8587133550488599
$ws.Range("F2").Value = 0.7071613301699083
$ws.Range("G2").Value = 0.7039249146757679
$ws.Range("H2").Value = 0.8519624573378839
